# Update "想去人数" (F column) figures across sheets to reflect the
# newly generated output (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 577
$ws1.Range("F5").Value = 133
$ws1.Range("F6").Value = 106
$ws1.Range("F7").Value = 279
$ws1.Range("F8").Value = 81
$ws1.Range("F10").Value = 16865
$ws1.Range("F14").Value = 6496
$ws1.Range("F21").Value = 97
$ws1.Range("F24").Value = 27
$ws1.Range("F25").Value = 15
$ws1.Range("F27").Value = 241
$ws1.Range("F28").Value = 920
$ws1.Range("F29").Value = 79
$ws1.Range("F30").Value = 5083
$ws1.Range("F32").Value = 45
$ws1.Range("F33").Value = 11551
$ws1.Range("F37").Value = 235
$ws1.Range("F38").Value = 3871

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 24

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 577
$ws4.Range("F5").Value = 133
$ws4.Range("F6").Value = 106
$ws4.Range("F7").Value = 279
$ws4.Range("F8").Value = 81
$ws4.Range("F10").Value = 16865
$ws4.Range("F14").Value = 6496
$ws4.Range("F21").Value = 97
$ws4.Range("F24").Value = 27
$ws4.Range("F25").Value = 15
$ws4.Range("F27").Value = 241
$ws4.Range("F28").Value = 920
$ws4.Range("F29").Value = 79
$ws4.Range("F30").Value = 5083
$ws4.Range("F32").Value = 24
$ws4.Range("F33").Value = 45
$ws4.Range("F34").Value = 11552
$ws4.Range("F38").Value = 235
$ws4.Range("F39").Value = 3871
